# Add a "2022-Q3" sheet (duplicate of the existing latest quarter sheet layout)
# right after the summary sheet, pushing "2022-Q2" .. "2020-Q4" down by one
# position, and fill in the new quarter's figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new top data row for
#    2022-Q3 and shift the existing quarterly rows down by one.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# New row 9 needs the same formatting as the existing index column (A2:A8).
$summary.Range("A8").Copy()
$summary.Range("A9").PasteSpecial(-4122)

# Shift the quarter/count/value columns down one row (bottom-up so we never
# clobber a row before it has been copied to its new location).
for ($r = 8; $r -ge 2; $r--) {
    $dst = $r + 1
    $summary.Range("B$dst").Value = $summary.Range("B$r").Value2
    $summary.Range("C$dst").Value = $summary.Range("C$r").Value2
    $summary.Range("D$dst").Value = $summary.Range("D$r").Value2
}

# Recompute the sequential index column (A2..A9 => 0..7).
for ($r = 2; $r -le 9; $r++) {
    $summary.Range("A$r").Value = $r - 2
}

# Fill in the brand new first row with the 2022-Q3 totals.
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.1

# ---------------------------------------------------------------------
# 2. Duplicate the current "2022-Q2" detail sheet so its data is preserved
#    under its own tab, then turn the original tab into "2022-Q3" and
#    overwrite it with the new quarter's fund holdings.
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item(2)

# Place the copy right after the current sheet; it keeps the "2022-Q2"
# data (and tab position 3) while the original (tab position 2) becomes
# the new quarter.
$q2Sheet.Copy($null, $q2Sheet)

$q3Sheet = $wb.Worksheets.Item(2)
$copySheet = $wb.Worksheets.Item(3)

# Avoid a transient name collision while renaming.
$q3Sheet.Name = "2022-Q3-tmp"
$copySheet.Name = "2022-Q2"
$q3Sheet.Name = "2022-Q3"

# Overwrite the (copied-from-Q2) values on the new 2022-Q3 sheet with the
# actual Q3 figures for the two funds already listed.
$q3Sheet.Range("D2").Value = 2.64
$q3Sheet.Range("E2").Value = "98.45"
$q3Sheet.Range("F2").Value = "2.82"
$q3Sheet.Range("G2").Value = "0.0744"

$q3Sheet.Range("C3").Value = "招商中证全指软件ETF"
$q3Sheet.Range("D3").Value = "1.02"
$q3Sheet.Range("E3").Value = "99.08"
$q3Sheet.Range("F3").Value = "2.85"
$q3Sheet.Range("G3").Value = "0.0291"
